# Slitrk2-Ptprs.xlsx: refresh with new TPM-derived NATMI numbers.
#
# Column D (Target cluster) keeps the same cluster names row-by-row; only
# the "Sending cluster" (column A) for the two 4-row blocks and the
# per-edge statistics (columns E..T) change to reflect the new TPM values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shared Ligand-side stats (columns E-J) for each 4-row sending-cluster block.
$ligandStatsByBlock = @{
    1 = @{ E = 2; F = 0.6666666666666666; G = 0.1735436666666667; H = 0.520631;  I = 0.4935485124749614; J = 0.4935485124749614 }
    2 = @{ E = 1; F = 0.3333333333333333; G = 0.1780806666666667; H = 0.534242;  I = 0.5064514875250385; J = 0.5064514875250385 }
}

# Per-row Receptor-side / edge stats (columns M-T). K and L are unchanged (3, 1).
$rowStats = @{
    2 = @{ M = 3.556762333333333;  N = 10.670287;        O = 0.04280930450251701; P = 0.04280930450251701; Q = 0.6172535767885555; R = 5.555282191097;    S = 0.02112846855730494; T = 0.02112846855730494 }
    3 = @{ M = 47.24901333333333;  N = 141.74704;         O = 0.5686906263805706;  P = 0.5686906263805704;  Q = 8.199767020248888;  R = 73.79790318223999; S = 0.2806764127085846;  T = 0.2806764127085846 }
    4 = @{ M = 24.53173066666666;  N = 73.595192;         O = 0.2952646900921413;  P = 0.2952646900921412;  Q = 4.257326489572444;  R = 38.31593840615199; S = 0.1457274485813568;  T = 0.1457274485813568 }
    5 = @{ M = 7.746355333333334;  N = 23.239066;         O = 0.09323537902477132; P = 0.0932353790247713;  Q = 1.344330907849556;  R = 12.098978170646;   S = 0.0460161826277151;  T = 0.04601618262771509 }
    6 = @{ M = 3.556762333333333;  N = 10.670287;         O = 0.04280930450251701; P = 0.04280930450251701; Q = 0.6333906074948888; R = 5.700515467454;    S = 0.02168083594521207; T = 0.02168083594521207 }
    7 = @{ M = 47.24901333333333;  N = 141.74704;         O = 0.5686906263805706;  P = 0.5686906263805704;  Q = 8.414135793742222;  R = 75.72722214368;    S = 0.2880142136719859;  T = 0.2880142136719858 }
    8 = @{ M = 24.53173066666666;  N = 73.595192;         O = 0.2952646900921413;  P = 0.2952646900921412;  Q = 4.36862695160711;   R = 39.31764256446399; S = 0.1495372415107845;  T = 0.1495372415107844 }
    9 = @{ M = 7.746355333333334;  N = 23.239066;         O = 0.09323537902477132; P = 0.1117885923419141;  Q = 1.379476121996889;  R = 12.415285097972;   S = 0.04721919639705621; T = 0.04721919639705621 }
}

# Rows 2-5: sending cluster "ECs"; rows 6-9: sending cluster "FAPs".
$sendingClusterByBlock = @{ 1 = "ECs"; 2 = "FAPs" }

for ($row = 2; $row -le 9; $row++) {
    $block = [int][Math]::Ceiling(($row - 1) / 4.0)

    $ws.Range("A$row").Value = $sendingClusterByBlock[$block]

    $ligand = $ligandStatsByBlock[$block]
    $ws.Range("E$row").Value = $ligand.E
    $ws.Range("F$row").Value = $ligand.F
    $ws.Range("G$row").Value = $ligand.G
    $ws.Range("H$row").Value = $ligand.H
    $ws.Range("I$row").Value = $ligand.I
    $ws.Range("J$row").Value = $ligand.J

    $stats = $rowStats[$row]
    $ws.Range("M$row").Value = $stats.M
    $ws.Range("N$row").Value = $stats.N
    $ws.Range("O$row").Value = $stats.O
    $ws.Range("P$row").Value = $stats.P
    $ws.Range("Q$row").Value = $stats.Q
    $ws.Range("R$row").Value = $stats.R
    $ws.Range("S$row").Value = $stats.S
    $ws.Range("T$row").Value = $stats.T
}
